$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1425304360311941
$ws.Range("D2").Value = 0.00328276805978156
$ws.Range("E2").Value = 0.4354855266417559
$ws.Range("F2").Value = 0.4268455756567278
$ws.Range("G2").Value = 0.3367752227119354
$ws.Range("H2").Value = 0.3434920082859492
$ws.Range("N2").Value = 1.995241137946778
$ws.Range("O2").Value = 1.288635611322178

$ws.Range("B3").Value = 0.1329091933101978
$ws.Range("D3").Value = 0.003117985810526136
$ws.Range("E3").Value = 0.3796236440278449
$ws.Range("F3").Value = 0.3933003553665202
$ws.Range("G3").Value = 0.3012008603001277
$ws.Range("H3").Value = 0.3306581927846537
$ws.Range("N3").Value = 1.846526605564804
$ws.Range("O3").Value = 1.186000399322808

$ws.Range("B4").Value = 0.1270661492628875
$ws.Range("D4").Value = 0.003018255499195988
$ws.Range("E4").Value = 0.3454298728940586
$ws.Range("F4").Value = 0.3729420886923549
$ws.Range("G4").Value = 0.2794993337179079
$ws.Range("H4").Value = 0.3229880749787242
$ws.Range("N4").Value = 1.755520437593759
$ws.Range("O4").Value = 1.123707417123967

$ws.Range("B5").Value = 0.1247015078381395
$ws.Range("D5").Value = 0.002977983059874845
$ws.Range("E5").Value = 0.3315199486717546
$ws.Range("F5").Value = 0.3647054950681934
$ws.Range("G5").Value = 0.2706907303260948
$ws.Range("H5").Value = 0.3199150623406126
$ws.Range("N5").Value = 1.718516884911963
$ws.Range("O5").Value = 1.098503372127993

$ws.Range("B6").Value = 0.1243098614095572
$ws.Range("D6").Value = 0.002971318247256605
$ws.Range("E6").Value = 0.3292116194363217
$ws.Range("F6").Value = 0.3633413997772692
$ws.Range("G6").Value = 0.269230162170615
$ws.Range("H6").Value = 0.3194079646665671
$ws.Range("N6").Value = 1.712377599108464
$ws.Range("O6").Value = 1.094329137056548

$ws.Range("B7").Value = 0.1270341919673399
$ws.Range("D7").Value = 0.003017710873152168
$ws.Range("E7").Value = 0.3452421832880503
$ws.Range("F7").Value = 0.3728307665163442
$ws.Range("G7").Value = 0.2793803970864275
$ws.Range("H7").Value = 0.3229464183999653
$ws.Range("N7").Value = 1.755021054648495
$ws.Range("O7").Value = 1.12336677637748

$ws.Range("B8").Value = 0.1391997998654944
$ws.Range("D8").Value = 0.00322565360309035
$ws.Range("E8").Value = 0.4162010468717057
$ws.Range("F8").Value = 0.4152293778963099
$ws.Range("G8").Value = 0.3244794999739042
$ws.Range("H8").Value = 0.3390232442816625
$ws.Range("N8").Value = 1.943904251618477
$ws.Range("O8").Value = 1.253095441307636

$ws.Range("B9").Value = 0.1635584291339143
$ws.Range("D9").Value = 0.003644731220816766
$ws.Range("E9").Value = 0.5562976445450403
$ws.Range("F9").Value = 0.5002909339772401
$ws.Range("G9").Value = 0.4140703990656505
$ws.Range("H9").Value = 0.3722245405283218
$ws.Range("N9").Value = 2.316504389646809
$ws.Range("O9").Value = 1.513334097579957

$ws.Range("B10").Value = 0.1817500012197257
$ws.Range("D10").Value = 0.003959326331720803
$ws.Range("E10").Value = 0.6599630079396519
$ws.Range("F10").Value = 0.5639956345128354
$ws.Range("G10").Value = 0.4806440831150098
$ws.Range("H10").Value = 0.3976543754028228
$ws.Range("N10").Value = 2.591327722520418
$ws.Range("O10").Value = 1.70823048577472

$ws.Range("B11").Value = 0.1900877858899008
$ws.Range("D11").Value = 0.004103863068408486
$ws.Range("E11").Value = 0.7073182102084701
$ws.Range("F11").Value = 0.5932485116357213
$ws.Range("G11").Value = 0.5111043360739984
$ws.Range("H11").Value = 0.4094518639923876
$ws.Range("N11").Value = 2.716530669508643
$ws.Range("O11").Value = 1.797728296575315

$ws.Range("B12").Value = 0.1932538308938092
$ws.Range("D12").Value = 0.004158796915575635
$ws.Range("E12").Value = 0.7252814645633805
$ws.Range("F12").Value = 0.604365721277091
$ws.Range("G12").Value = 0.5226648403903766
$ws.Range("H12").Value = 0.4139524800987431
$ws.Range("N12").Value = 2.763963193734583
$ws.Range("O12").Value = 1.831741466706262

$ws.Range("B13").Value = 0.1925715833279469
$ws.Range("D13").Value = 0.004146957056267553
$ws.Range("E13").Value = 0.7214113434755376
$ws.Range("F13").Value = 0.6019696543234119
$ws.Range("G13").Value = 0.5201739184791165
$ws.Range("H13").Value = 0.4129817139520924
$ws.Range("N13").Value = 2.753746911385861
$ws.Range("O13").Value = 1.824410655378244

$ws.Range("B14").Value = 0.190348085298524
$ws.Range("D14").Value = 0.004108378505360832
$ws.Range("E14").Value = 0.7087954249356869
$ws.Range("F14").Value = 0.5941623311389179
$ws.Range("G14").Value = 0.5120549037030742
$ws.Range("H14").Value = 0.409821466670877
$ws.Range("N14").Value = 2.720432585385993
$ws.Range("O14").Value = 1.800524119720649

$ws.Range("B15").Value = 0.1889872547792919
$ws.Range("D15").Value = 0.00408477407517438
$ws.Range("E15").Value = 0.7010719049055893
$ws.Range("F15").Value = 0.5893853151532795
$ws.Range("G15").Value = 0.5070851600557091
$ws.Range("H15").Value = 0.4078900483665961
$ws.Range("N15").Value = 2.700029152110062
$ws.Range("O15").Value = 1.785908904299276

$ws.Range("B16").Value = 0.1812063398625838
$ws.Range("D16").Value = 0.003949908852845851
$ws.Range("E16").Value = 0.6568724169760856
$ws.Range("F16").Value = 0.5620894304974513
$ws.Range("G16").Value = 0.4786570196494324
$ws.Range("H16").Value = 0.3968880122874339
$ws.Range("N16").Value = 2.583148644061055
$ws.Range("O16").Value = 1.702398612130764

$ws.Range("B17").Value = 0.1764487851811651
$ws.Range("D17").Value = 0.003867535620365459
$ws.Range("E17").Value = 0.629809763319642
$ws.Range("F17").Value = 0.5454146261146775
$ws.Range("G17").Value = 0.4612626808594769
$ws.Range("H17").Value = 0.3901974874650023
$ws.Range("N17").Value = 2.511489453097226
$ws.Range("O17").Value = 1.651383722660285

$ws.Range("B18").Value = 0.1737182494248515
$ws.Range("D18").Value = 0.003820291157840217
$ws.Range("E18").Value = 0.6142624613427046
$ws.Range("F18").Value = 0.5358494138644119
$ws.Range("G18").Value = 0.4512744032264209
$ws.Range("H18").Value = 0.3863708487733675
$ws.Range("N18").Value = 2.470290675909951
$ws.Range("O18").Value = 1.622120123892898

$ws.Range("B19").Value = 0.1727947553757332
$ws.Range("D19").Value = 0.003804318212985436
$ws.Range("E19").Value = 0.6090015003013463
$ws.Range("F19").Value = 0.5326151957723084
$ws.Range("G19").Value = 0.4478953566152768
$ws.Range("H19").Value = 0.3850789144428575
$ws.Range("N19").Value = 2.456344689235948
$ws.Range("O19").Value = 1.612225456820568

$ws.Range("B20").Value = 0.1769546279153928
$ws.Range("D20").Value = 0.003876290507584201
$ws.Range("E20").Value = 0.6326887045236589
$ws.Range("F20").Value = 0.5471870247896504
$ws.Range("G20").Value = 0.4631126271809762
$ws.Range("H20").Value = 0.3909074710753089
$ws.Range("N20").Value = 2.519115905986666
$ws.Range("O20").Value = 1.656806175146926

$ws.Range("B21").Value = 0.1910009463247917
$ws.Range("D21").Value = 0.004119704540567426
$ws.Range("E21").Value = 0.7125001682736638
$ws.Range("F21").Value = 0.5964544477321994
$ws.Range("G21").Value = 0.5144389481118878
$ws.Range("H21").Value = 0.4107488063632729
$ws.Range("N21").Value = 2.730217288371648
$ws.Range("O21").Value = 1.807536841490844

$ws.Range("B22").Value = 0.200231658042739
$ws.Range("D22").Value = 0.00427995902020939
$ws.Range("E22").Value = 0.7648426698322197
$ws.Range("F22").Value = 0.6288856518111032
$ws.Range("G22").Value = 0.5481347357242612
$ws.Range("H22").Value = 0.4239096552953754
$ws.Range("N22").Value = 2.868303394937072
$ws.Range("O22").Value = 1.906761666784462

$ws.Range("B23").Value = 0.1953005061487829
$ws.Range("D23").Value = 0.00419432250900087
$ws.Range("E23").Value = 0.7368890890760866
$ws.Range("F23").Value = 0.6115551107154431
$ws.Range("G23").Value = 0.5301366292087835
$ws.Range("H23").Value = 0.4168677050772089
$ws.Range("N23").Value = 2.794595190495954
$ws.Range("O23").Value = 1.853737639956933

$ws.Range("B24").Value = 0.1767259219250121
$ws.Range("D24").Value = 0.003872332071402695
$ws.Range("E24").Value = 0.6313870999754982
$ws.Range("F24").Value = 0.5463856568671446
$ws.Range("G24").Value = 0.4622762292460152
$ws.Range("H24").Value = 0.3905864257188227
$ws.Range("N24").Value = 2.515667989176677
$ws.Range("O24").Value = 1.654354480355494

$ws.Range("B25").Value = 0.1569161479392847
$ws.Range("D25").Value = 0.003530170004244937
$ws.Range("E25").Value = 0.5182800985496812
$ws.Range("F25").Value = 0.4770700030462791
$ws.Range("G25").Value = 0.3897053291205168
$ws.Range("H25").Value = 0.3630617902832967
$ws.Range("N25").Value = 2.215498394811846
$ws.Range("O25").Value = 1.442293097308493
